$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 685; everything from 685 onward
# shifts down by 3 (685->688 ... 733->736), matching the diff's dimension
# change from A1:T733 to A1:T736.
$ws.Rows("685:687").Insert()

# --- New row 685 ---
$ws.Range("A685").Value = 11
$ws.Range("B685").Value = "Vega Monumental Concepción"
$ws.Range("C685").Value = "Bíobío"
$ws.Range("D685").Value = 45021
$ws.Range("E685").Value = 8
$ws.Range("F685").Value = "Fruta"
$ws.Range("G685").Value = 100102
$ws.Range("H685").Value = "Cítricos"
$ws.Range("I685").Value = 100102003
$ws.Range("J685").Value = "Limón"
$ws.Range("K685").Value = "Sin especificar"
$ws.Range("L685").Value = "1a plateado"
$ws.Range("M685").Value = 270
$ws.Range("N685").Value = 21000
$ws.Range("O685").Value = 22000
$ws.Range("P685").Value = 21444
$ws.Range("Q685").Value = "$/malla 16 kilos"
$ws.Range("R685").Value = "Provincia de Limarí"
$ws.Range("S685").Value = 1340
$ws.Range("T685").Value = 16

# --- New row 686 ---
$ws.Range("A686").Value = 11
$ws.Range("B686").Value = "Vega Monumental Concepción"
$ws.Range("C686").Value = "Bíobío"
$ws.Range("D686").Value = 45021
$ws.Range("E686").Value = 8
$ws.Range("F686").Value = "Fruta"
$ws.Range("G686").Value = 100102
$ws.Range("H686").Value = "Cítricos"
$ws.Range("I686").Value = 100102003
$ws.Range("J686").Value = "Limón"
$ws.Range("K686").Value = "Sin especificar"
$ws.Range("L686").Value = "1a plateado"
$ws.Range("M686").Value = 270
$ws.Range("N686").Value = 21000
$ws.Range("O686").Value = 22000
$ws.Range("P686").Value = 21556
$ws.Range("Q686").Value = "$/malla 16 kilos"
$ws.Range("R686").Value = "Región de O'Higgins"
$ws.Range("S686").Value = 1347
$ws.Range("T686").Value = 16

# --- New row 687 ---
$ws.Range("A687").Value = 11
$ws.Range("B687").Value = "Vega Monumental Concepción"
$ws.Range("C687").Value = "Bíobío"
$ws.Range("D687").Value = 45021
$ws.Range("E687").Value = 8
$ws.Range("F687").Value = "Fruta"
$ws.Range("G687").Value = 100102
$ws.Range("H687").Value = "Cítricos"
$ws.Range("I687").Value = 100102003
$ws.Range("J687").Value = "Limón"
$ws.Range("K687").Value = "Sin especificar"
$ws.Range("L687").Value = "2a plateado"
$ws.Range("M687").Value = 150
$ws.Range("N687").Value = 17000
$ws.Range("O687").Value = 17000
$ws.Range("P687").Value = 17000
$ws.Range("Q687").Value = "$/malla 16 kilos"
$ws.Range("R687").Value = "Región de O'Higgins"
$ws.Range("S687").Value = 1062
$ws.Range("T687").Value = 16
